# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-11, as recalculated after switching
# the streak metric from Strike# to K.
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
